$wb = $excel.ActiveWorkbook

# Sheet: Escapement
$ws = $wb.Worksheets.Item("Escapement")
$ws.Range("A60").Value = 4687888.21634433
$ws.Range("B60").Value = 3229999.83754533
$ws.Range("C60").Value = 2783049.2461762
$ws.Range("D60").Value = 1846975.3373724
$ws.Range("E60").Value = 2903009.95933991

# Sheet: Total Catch
$ws = $wb.Worksheets.Item("Total Catch")
$ws.Range("A60").Value = 4182270.2673438
$ws.Range("B60").Value = 2832609.70661295
$ws.Range("C60").Value = 4299080.64742782
$ws.Range("D60").Value = 5893351.00034585
$ws.Range("E60").Value = 7797211.51975421

# Sheet: Run Size
$ws = $wb.Worksheets.Item("Run Size")
$ws.Range("A60").Value = 8885790.26768812
$ws.Range("B60").Value = 6069513.70615829
$ws.Range("C60").Value = 7095614.64760402
$ws.Range("D60").Value = 7725547.00071825
$ws.Range("E60").Value = 10657141.5200941

# Sheet: Run Size no Offshore
$ws = $wb.Worksheets.Item("Run Size no Offshore")
$ws.Range("A60").Value = 8461494.92707114
$ws.Range("B60").Value = 5776402.50584379
$ws.Range("C60").Value = 6785280.23768978
$ws.Range("D60").Value = 7322840.46632597
$ws.Range("E60").Value = 10100440.3009972
